$wb = $excel.ActiveWorkbook
$wsSteps = $wb.Worksheets.Item("TestSteps")
$wsData  = $wb.Worksheets.Item("TestData")

# Insert a new test step row (row 3) in the TestSteps sheet, pushing the
# existing Keyword/ObjectID/KeyInData (columns A:C) rows down by one while
# leaving the KeyInData/screenshot/log columns (D:F) untouched - this mirrors
# what Excel does for "Insert Cells... Shift cells down" over A3:C3.
$wsSteps.Range("A3:C3").Insert(-4121) | Out-Null

# The insert leaves row 3 (A:F) with a freshly-minted, unformatted style;
# restore it by copying the formatting from row 4 (which now holds the
# original row-3 content/format) back onto row 3.
$wsSteps.Range("A4:F4").Copy()
$wsSteps.Range("A3:F3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# The partial-column insert also leaves a stray, fully-empty row 17 behind
# (the sheet only really uses rows 1-16) - remove it so the dimension stays
# A1:F16.
$wsSteps.Rows(17).Delete() | Out-Null

# Populate the newly inserted step: verify the "last processed date" field
# is not present on the page.
$wsSteps.Range("A3").Value = "verifyElementNotPresent"
$wsSteps.Range("B3").Value = "txt_last_proc_date"
$wsSteps.Range("C3").Value = "getData=WaitForPageLoad"

# Make TestSteps the active sheet/tab, with A15 as the selected cell
# (previously TestData was active, and B8 was selected on TestSteps).
$wsSteps.Select() | Out-Null
$wsSteps.Range("A15").Select() | Out-Null

# TestData keeps its own last selection (E2); it is simply no longer the
# active/visible tab.
$wsData.Range("E2").Select() | Out-Null
$wsSteps.Select() | Out-Null
